# CollMeet.xlsx — add new field "MeetPsnNote" (面晤人備註) to the DB layout table.
#
# The table lists DB columns (SEQ / 欄位名稱 / 中文名稱 / 形態 / 長度 / 小數 / 備註說明)
# for the CollMeet table, one row per field, starting at worksheet row 9.
# A new field MeetPsnNote / 面晤人備註 (NVARCHAR2, length 50) is inserted right
# after "Remark" (row 23) and before "CreateDate" (old row 24), i.e. as new
# row 24. Every row from the old row 24 onward shifts down by one, and the
# SEQ numbers of those shifted rows increase by one to stay sequential.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new blank row at 24 (old rows 24-27 shift down to 25-28) ---
$ws.Rows.Item(24).Insert()

# Copy the formatting (font, borders, alignment, number format) of the row
# above (row 23, "Remark") down onto the freshly inserted blank row so it
# matches the rest of the table.
$ws.Range("A23:G23").Copy()
$ws.Range("A24:G24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Fill in the new row's values ---
$ws.Range("A24").Value = 16
$ws.Range("B24").Value = "MeetPsnNote"
$ws.Range("C24").Value = "面晤人備註"
$ws.Range("D24").Value = "NVARCHAR2"
$ws.Range("E24").Value = 50
$ws.Range("G24").Value = "2023/8/23新增,舊資料由RECEIVE_PERSONNAME轉入"

# Highlight the new field name / Chinese name cells in yellow so the addition
# stands out in the layout sheet.
$ws.Range("B24:C24").Interior.Color = 65535

# --- 3. Renumber the SEQ column for the rows that shifted down ---
$ws.Range("A25").Value = 17
$ws.Range("A26").Value = 18
$ws.Range("A27").Value = 19
$ws.Range("A28").Value = 20

# --- 4. Update the view: selection moves to the newly added cell, scrolled
#        so the new row is visible near the top of the window ---
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("B24").Select()
